$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.179.95"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "3.528.06"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.11"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.20"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("D7").Value = "3.527.76"
$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.85"
$ws.Range("E11").Value = "  +2.74%  "

$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "4.120.21"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.04"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "3.517.36"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "67.235.32"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").Value = "  +8.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.44"
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.40"
$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.32"
$ws.Range("E22").Value = "  -2.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.611"
$ws.Range("E23").Value = "  -2.76%  "

$ws.Range("E24").Value = "  +1.66%  "

$ws.Range("D25").Value = "3.663.37"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  -4.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.36"
$ws.Range("E29").Value = "  -4.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.59"
$ws.Range("E31").Value = "  -3.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  -2.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.47"
$ws.Range("E34").Value = "  -0.48%  "

$ws.Range("D35").Value = "3.520.67"
$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("E36").Value = "  -2.98%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.95"
$ws.Range("E37").Value = "  -3.26%  "

$ws.Range("E38").Value = "  +1.38%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0897"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "171.87"
$ws.Range("E42").Value = "  -4.01%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("E44").Value = "  -9.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.898"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.07"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.51"
$ws.Range("E47").Value = "  -4.81%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.28"
$ws.Range("E48").Value = "  -1.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.50"
$ws.Range("E49").Value = "  -1.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.46"
$ws.Range("E50").Value = "  -4.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.993"
$ws.Range("E51").Value = "  +0.23%  "
